$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 190942.08
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 196232.14
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 588696.42
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -589032.42

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4463.0386
$ws.Range("J112").Value = 4729.3184
$ws.Range("L112").Value = 14187.9552
$ws.Range("N112").Value = -16403.9552

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1322159.5
$ws.Range("I137").Value = 6553.3335
$ws.Range("J137").Value = 2180163.5
$ws.Range("K137").Value = 19660.0005
$ws.Range("L137").Value = 6540490.5
$ws.Range("M137").Value = -17110.0005
$ws.Range("N137").Value = -6545590.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 18174.375
$ws.Range("J138").Value = 4565.8335
$ws.Range("L138").Value = 13697.5005
$ws.Range("N138").Value = -23977.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 604383.9
$ws.Range("I74").Value = 3343.1875
$ws.Range("J74").Value = 4451044
$ws.Range("K74").Value = 3343.1875
$ws.Range("L74").Value = 4451044
$ws.Range("M74").Value = -2469.1875
$ws.Range("N74").Value = -4452792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 604383.9
$ws.Range("I77").Value = 3343.1875
$ws.Range("J77").Value = 4451044
$ws.Range("K77").Value = 16715.9375
$ws.Range("L77").Value = 22255220
$ws.Range("M77").Value = -12347.9375
$ws.Range("N77").Value = -22263956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H105").Value = 95923.08
$ws.Range("I105").Value = 95923.08
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 95923.08
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -92429.08
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8072
$ws.Range("I105").Value = 5587.2593
$ws.Range("K105").Value = 5587.2593
$ws.Range("M105").Value = -3840.2593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 595.6667
$ws.Range("I22").Value = 609.6957
$ws.Range("K22").Value = 609.6957
$ws.Range("M22").Value = -259.6957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4431.9897
$ws.Range("I31").Value = 2694.9048
$ws.Range("J31").Value = 4899.6665
$ws.Range("K31").Value = 2694.9048
$ws.Range("L31").Value = 4899.6665
$ws.Range("M31").Value = -2399.9048
$ws.Range("N31").Value = -5489.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4431.9897
$ws.Range("I34").Value = 2694.9048
$ws.Range("J34").Value = 4899.6665
$ws.Range("K34").Value = 2694.9048
$ws.Range("L34").Value = 4899.6665
$ws.Range("M34").Value = -2492.9048
$ws.Range("N34").Value = -5303.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1252653
$ws.Range("I99").Value = 1820751.2
$ws.Range("J99").Value = 2836.8
$ws.Range("K99").Value = 1820751.2
$ws.Range("L99").Value = 2836.8
$ws.Range("M99").Value = -1819253.2
$ws.Range("N99").Value = -5832.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1252653
$ws.Range("I126").Value = 1820751.2
$ws.Range("J126").Value = 2836.8
$ws.Range("K126").Value = 5462253.6
$ws.Range("L126").Value = 8510.400000000001
$ws.Range("M126").Value = -5459783.6
$ws.Range("N126").Value = -13450.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2869.8333
$ws.Range("I132").Value = 2933.6667
$ws.Range("K132").Value = 8801.000100000001
$ws.Range("M132").Value = -6271.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2836.9
$ws.Range("I50").Value = 3400.2
$ws.Range("J50").Value = 2273.6
$ws.Range("K50").Value = 10200.6
$ws.Range("L50").Value = 6820.799999999999
$ws.Range("M50").Value = -9719.599999999999
$ws.Range("N50").Value = -7782.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 2836.9
$ws.Range("I53").Value = 3400.2
$ws.Range("J53").Value = 2273.6
$ws.Range("K53").Value = 10200.6
$ws.Range("L53").Value = 6820.799999999999
$ws.Range("M53").Value = -9719.599999999999
$ws.Range("N53").Value = -7782.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2498.9285
$ws.Range("J68").Value = 3328.3333
$ws.Range("L68").Value = 9984.999899999999
$ws.Range("N68").Value = -11606.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 17979.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 17979.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 53938.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -54568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2498.9285
$ws.Range("J71").Value = 3328.3333
$ws.Range("L71").Value = 29954.9997
$ws.Range("N71").Value = -38066.9997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 17979.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 17979.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 53938.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -56122.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2312.8
$ws.Range("I75").Value = 952
$ws.Range("J75").Value = 2896
$ws.Range("K75").Value = 2856
$ws.Range("L75").Value = 8688
$ws.Range("M75").Value = -1858
$ws.Range("N75").Value = -10684

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2312.8
$ws.Range("I78").Value = 952
$ws.Range("J78").Value = 2896
$ws.Range("K78").Value = 8568
$ws.Range("L78").Value = 26064
$ws.Range("M78").Value = -3576
$ws.Range("N78").Value = -36048

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 567.6667
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 551.5
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 1654.5
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -4150.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 6821.4585
$ws.Range("I129").Value = 1466.8889
$ws.Range("J129").Value = 10034.2
$ws.Range("K129").Value = 4400.6667
$ws.Range("L129").Value = 30102.6
$ws.Range("M129").Value = 599.3333000000002
$ws.Range("N129").Value = -40102.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3789543.5
$ws.Range("J131").Value = 1923.2307
$ws.Range("L131").Value = 5769.6921
$ws.Range("N131").Value = -15849.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 53356
$ws.Range("J105").Value = 53356
$ws.Range("L105").Value = 53356
$ws.Range("N105").Value = -60344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 60000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -62524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 10210000
$ws.Range("I20").Value = 212499.75
$ws.Range("J20").Value = 50200000
$ws.Range("K20").Value = 212499.75
$ws.Range("L20").Value = 50200000
$ws.Range("M20").Value = -212273.75
$ws.Range("N20").Value = -50200452

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6742.0835
$ws.Range("I22").Value = 813
$ws.Range("J22").Value = 7927.9
$ws.Range("K22").Value = 813
$ws.Range("L22").Value = 7927.9
$ws.Range("M22").Value = -518
$ws.Range("N22").Value = -8517.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 6742.0835
$ws.Range("I27").Value = 813
$ws.Range("J27").Value = 7927.9
$ws.Range("K27").Value = 813
$ws.Range("L27").Value = 7927.9
$ws.Range("M27").Value = -706
$ws.Range("N27").Value = -8141.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 662.75757
$ws.Range("I55").Value = 484.1875
$ws.Range("K55").Value = 484.1875
$ws.Range("M55").Value = -311.1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 74999.664
$ws.Range("J96").Value = 74999.664
$ws.Range("L96").Value = 74999.664
$ws.Range("N96").Value = -80491.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2931.625
$ws.Range("I100").Value = 2881.889
$ws.Range("J100").Value = 2995.5715
$ws.Range("K100").Value = 2881.889
$ws.Range("L100").Value = 2995.5715
$ws.Range("M100").Value = -2340.889
$ws.Range("N100").Value = -4077.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 33122.5
$ws.Range("J101").Value = 33122.5
$ws.Range("L101").Value = 33122.5
$ws.Range("N101").Value = -39612.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 31065.8
$ws.Range("J103").Value = 31065.8
$ws.Range("L103").Value = 31065.8
$ws.Range("N103").Value = -33409.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 5000
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 5000
$ws.Range("I24").Value = 5000
$ws.Range("K24").Value = 5000
$ws.Range("M24").Value = -4770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 7522500
$ws.Range("I53").Value = 7522500
$ws.Range("K53").Value = 7522500
$ws.Range("M53").Value = -7521893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1981.6666
$ws.Range("I62").Value = 1981.6666
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1981.6666
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1357.6666
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 1981.6666
$ws.Range("I65").Value = 1981.6666
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9908.333000000001
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6788.333000000001
$ws.Range("N65").ClearContents()
